# Applies the "№" (row number / index) column to the "guns" and "armor"
# worksheets, and updates the active-sheet / selection bookkeeping to match
# the authored diff (guns becomes the active tab with a selection of F14,
# armor's selection moves to A7).

$wb = $excel.ActiveWorkbook

$header = [char]8470   # "№" (U+2116 NUMERO SIGN)

# ---------------------------------------------------------------------
# "guns" sheet (sheet2.xml): add column F ("№") with a 0-based row index
# ---------------------------------------------------------------------
$wsGuns = $wb.Worksheets.Item("guns")

$wsGuns.Range("F1").Value = $header
for ($row = 2; $row -le 13; $row++) {
    $wsGuns.Cells.Item($row, 6).Value = $row - 2
}

# ---------------------------------------------------------------------
# "armor" sheet (sheet4.xml): add column F ("№") with a 0-based row index
# ---------------------------------------------------------------------
$wsArmor = $wb.Worksheets.Item("armor")

$wsArmor.Range("F1").Value = $header
for ($row = 2; $row -le 12; $row++) {
    $wsArmor.Cells.Item($row, 6).Value = $row - 2
}

# ---------------------------------------------------------------------
# Selection / active-tab bookkeeping
# ---------------------------------------------------------------------
# "armor" keeps its own (non-active) selection, now on A7.
$wsArmor.Range("A7").Select() | Out-Null

# "guns" becomes the active sheet/tab, selection on F14 (this also clears
# tabSelected on whichever sheet previously had it, i.e. "enemies").
$wsGuns.Range("F14").Select() | Out-Null
